# edit.ps1 - reproduces the "Various updates needing MeOH investigation" commit
# against Paper/excel figures.xlsx

$wb = $excel.ActiveWorkbook

$wsGrowth    = $wb.Worksheets.Item("Growth Yield")
$wsKO        = $wb.Worksheets.Item("Knockout Validation")
$wsStats     = $wb.Worksheets.Item("General Model Stats")
$wsMetab     = $wb.Worksheets.Item("Metabolomics Figure")

# ---------------------------------------------------------------------------
# 1. "Growth Yield" sheet - add two new summary rows (Model Error / Model Top
#    Error) below the existing Predicted/Experimental/+-Error rows.
# ---------------------------------------------------------------------------

$wsGrowth.Range("A8").Value = "Model Error"
$wsGrowth.Range("B8").Formula = "=(B3-B4)/B4"
$wsGrowth.Range("C8").Formula = "=(C3-C4)/C4"
$wsGrowth.Range("B8:C8").NumberFormat = "0%"

$wsGrowth.Range("A9").Value = "Model Top Error"
$wsGrowth.Range("B9").Formula = "=(B3-B6)/B6"
$wsGrowth.Range("C9").Formula = "=(C3-C6)/C6"
$wsGrowth.Range("B9:C9").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# 2. "Knockout Validation" sheet - move the TOTAL figure down one row and add
#    a new MCC (Matthews correlation coefficient) figure in its old spot.
# ---------------------------------------------------------------------------

$wsKO.Range("G17").Value = "MCC"
$wsKO.Range("G17").Font.Bold = $false
$wsKO.Range("H17").Value = 0.67

$wsKO.Range("G19").Value = "TOTAL"

# ---------------------------------------------------------------------------
# 3. "General Model Stats" sheet - refreshed model (iMR494 -> iMR524) numbers.
# ---------------------------------------------------------------------------

$wsStats.Range("C4").Value = 524
$wsStats.Range("C6").Value = "646/45"
$wsStats.Range("C7").Value = 268
$wsStats.Range("C8").Value = 545
$wsStats.Range("C9").Value = 53
$wsStats.Range("C10").Value = 489
$wsStats.Range("B12").Value = "Table 1. A comparison between iMR524 and iMM518"

# ---------------------------------------------------------------------------
# 4. Selections / active sheet bookkeeping.
# ---------------------------------------------------------------------------

$wsKO.Range("H18").Select()

$wsStats.Range("B12").Select()

$wsGrowth.Activate()
$wsGrowth.Range("D13").Select()
